# Add two new columns (I: "I0", J: "IF") to the header row and their
# corresponding data-row values (8, 8), matching the style already used
# by the existing header cells (bold, centered, bordered).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header cell H1 (which already carries the header
# style) into I1 and J1, then overwrite their values/text so the style
# comes along for the ride.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("I1").Value = "I0"

$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("J1").Value = "IF"

# New data-row values under the new headers.
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8
